# Menu-Languages.docx edit: add several new menu-item lines and rename
# one existing line ("Graphics" -> "Marketing Materials"), per the
# supplied unified diff.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-ParagraphAfterText {
    # Locates the paragraph whose text exactly equals $anchorText (after
    # trimming the trailing paragraph mark), inserts a brand-new empty
    # paragraph right after it, then stamps that new paragraph with the
    # supplied raw <w:p>...</w:p> OOXML via Range.InsertXML.
    param(
        [string]$anchorText,
        [string]$paraXml
    )

    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Insert-ParagraphAfterText: anchor not found: $anchorText"
    }
    $anchorPara = $searchRange.Paragraphs.First
    $anchorIndex = $anchorPara.Index

    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($anchorIndex + 1)
    $newRange = $newPara.Range
    $newRange.InsertXML($pkgHeader + $paraXml + $pkgFooter)
}

# --- 1. New "The Other Side" line, right after the "SmartCard" entry
#        that precedes "Resources" -------------------------------------
$theOtherSideXml = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">The Other Side</w:t></w:r></w:p>
'@

Insert-ParagraphAfterText "SmartCard" $theOtherSideXml

# --- 2. The pre-existing "Graphics" line (indent left=360 firstLine=360)
#        is retitled "Marketing Materials" and gains a leading tab.
#        Do this BEFORE any new paragraphs are inserted nearby, while
#        "Graphics" is still a unique anchor in the document. ------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Graphics", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find original 'Graphics' paragraph to rename"
}
$graphicsPara = $searchRange.Paragraphs.First
$graphicsRunRange = $graphicsPara.Range
$marketingMaterialsRunXml = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">Marketing Materials</w:t></w:r></w:p>'
$graphicsRunRange.InsertXML($pkgHeader + $marketingMaterialsRunXml + $pkgFooter)

# --- 3. New "Exchanges Listing Guide" line, right after the "Exchanges"
#        entry ------------------------------------------------------------
$exchangesListingGuideXml = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:tab/><w:t xml:space="preserve">Exchanges Listing Guide</w:t></w:r></w:p>
'@

$tab = [char]9
Insert-ParagraphAfterText ($tab + "Exchanges") $exchangesListingGuideXml

# --- 4. New "Graphics" line, right after the just-inserted
#        "Exchanges Listing Guide" line (lands right before the renamed
#        "Marketing Materials" paragraph, as required) --------------------
$graphicsXml = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Graphics</w:t></w:r></w:p>
'@

Insert-ParagraphAfterText ("$tab$tab" + "Exchanges Listing Guide") $graphicsXml

# --- 5. New "Guide" line, right after the "Electrum Wallet" entry --------
$guideXml = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:tab/><w:t>Guide</w:t></w:r></w:p>
'@

Insert-ParagraphAfterText "Electrum Wallet" $guideXml

Write-Host "Done."
